$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 5649.7236
$ws.Range("I32").Value = 4684.1333
$ws.Range("J32").Value = 7353.706
$ws.Range("K32").Value = 4684.1333
$ws.Range("L32").Value = 7353.706
$ws.Range("M32").Value = -4397.1333
$ws.Range("N32").Value = -7927.706

# Row 44
$ws.Range("H44").Value = 19996.666
$ws.Range("J44").Value = 19996.666
$ws.Range("L44").Value = 19996.666
$ws.Range("N44").Value = -20972.666

# Row 55
$ws.Range("H55").Value = 21586.334
$ws.Range("I55").Value = 0
$ws.Range("J55").Value = 21586.334
$ws.Range("K55").Value = 0
$ws.Range("L55").Value = 21586.334
$ws.Range("M55").ClearContents()
$ws.Range("N55").Value = -22216.334

# Row 74
$ws.Range("H74").Value = 4542.9473
$ws.Range("I74").Value = 5013.2583
$ws.Range("J74").Value = 2460.1428
$ws.Range("K74").Value = 5013.2583
$ws.Range("L74").Value = 2460.1428
$ws.Range("M74").Value = -4139.2583
$ws.Range("N74").Value = -4208.1428

# Row 77
$ws.Range("H77").Value = 4542.9473
$ws.Range("I77").Value = 5013.2583
$ws.Range("J77").Value = 2460.1428
$ws.Range("K77").Value = 25066.2915
$ws.Range("L77").Value = 12300.714
$ws.Range("M77").Value = -20698.2915
$ws.Range("N77").Value = -21036.714

# Row 97
$ws.Range("H97").Value = 787.1429000000001
$ws.Range("I97").Value = 527.25
$ws.Range("J97").Value = 1133.6666
$ws.Range("K97").Value = 527.25
$ws.Range("L97").Value = 1133.6666
$ws.Range("M97").Value = -31.25
$ws.Range("N97").Value = -2125.6666

# Row 110
$ws.Range("H110").Value = 2000
$ws.Range("I110").Value = 2000
$ws.Range("K110").Value = 2000
$ws.Range("M110").Value = 45

# Row 122
$ws.Range("H122").Value = 2138486.5
$ws.Range("I122").Value = 2850693
$ws.Range("J122").Value = 1866.6666
$ws.Range("K122").Value = 8552079
$ws.Range("L122").Value = 5599.9998
$ws.Range("M122").Value = -8549629
$ws.Range("N122").Value = -10499.9998

# Row 132
$ws.Range("H132").Value = 5189.9287
$ws.Range("I132").Value = 1506.0646
$ws.Range("J132").Value = 15571.728
$ws.Range("K132").Value = 4518.1938
$ws.Range("L132").Value = 46715.18399999999
$ws.Range("M132").Value = -1988.1938
$ws.Range("N132").Value = -51775.18399999999

# Row 139
$ws.Range("H139").Value = 44000
$ws.Range("J139").Value = 44000
$ws.Range("L139").Value = 44000
$ws.Range("N139").Value = -54280

$ws = $wb.Worksheets.Item("BSM")
# Row 107
$ws.Range("H107").Value = 1741.1904
$ws.Range("I107").Value = 1645.6666
$ws.Range("J107").Value = 1980
$ws.Range("K107").Value = 1645.6666
$ws.Range("L107").Value = 1980
$ws.Range("M107").Value = 274.3334
$ws.Range("N107").Value = -5820

# Row 134
$ws.Range("H134").Value = 6712.923
$ws.Range("I134").Value = 11188
$ws.Range("J134").Value = 2877.1428
$ws.Range("K134").Value = 33564
$ws.Range("L134").Value = 8631.428400000001
$ws.Range("M134").Value = -31029
$ws.Range("N134").Value = -13701.4284

$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 6426.8696
$ws.Range("I31").Value = 1547.2941
$ws.Range("J31").Value = 20252.334
$ws.Range("K31").Value = 1547.2941
$ws.Range("L31").Value = 20252.334
$ws.Range("M31").Value = -1252.2941
$ws.Range("N31").Value = -20842.334

# Row 34
$ws.Range("H34").Value = 6426.8696
$ws.Range("I34").Value = 1547.2941
$ws.Range("J34").Value = 20252.334
$ws.Range("K34").Value = 1547.2941
$ws.Range("L34").Value = 20252.334
$ws.Range("M34").Value = -1345.2941
$ws.Range("N34").Value = -20656.334

# Row 107
$ws.Range("H107").Value = 937.7895
$ws.Range("I107").Value = 751.2857
$ws.Range("J107").Value = 1046.5834
$ws.Range("K107").Value = 751.2857
$ws.Range("L107").Value = 1046.5834
$ws.Range("M107").Value = 1168.7143
$ws.Range("N107").Value = -4886.5834

# Row 122
$ws.Range("H122").Value = 954.3158
$ws.Range("I122").Value = 902.6667
$ws.Range("J122").Value = 1000.8
$ws.Range("K122").Value = 2708.0001
$ws.Range("L122").Value = 3002.4
$ws.Range("M122").Value = -258.0001000000002
$ws.Range("N122").Value = -7902.4

$ws = $wb.Worksheets.Item("CUL")
# Row 25
$ws.Range("H25").Value = 3183.5
$ws.Range("I25").Value = 399.5
$ws.Range("J25").Value = 4575.5
$ws.Range("K25").Value = 1198.5
$ws.Range("L25").Value = 13726.5
$ws.Range("M25").Value = -1029.5
$ws.Range("N25").Value = -14064.5

# Row 30
$ws.Range("H30").Value = 3183.5
$ws.Range("I30").Value = 399.5
$ws.Range("J30").Value = 4575.5
$ws.Range("K30").Value = 1198.5
$ws.Range("L30").Value = 13726.5
$ws.Range("M30").Value = -1096.5
$ws.Range("N30").Value = -13930.5

# Row 62
$ws.Range("H62").Value = 5250
$ws.Range("I62").Value = 4000
$ws.Range("J62").Value = 6500
$ws.Range("K62").Value = 12000
$ws.Range("L62").Value = 19500
$ws.Range("M62").Value = -11314
$ws.Range("N62").Value = -20872

# Row 65
$ws.Range("H65").Value = 5250
$ws.Range("I65").Value = 4000
$ws.Range("J65").Value = 6500
$ws.Range("K65").Value = 36000
$ws.Range("L65").Value = 58500
$ws.Range("M65").Value = -32568
$ws.Range("N65").Value = -65364

# Row 140
$ws.Range("H140").Value = 6880.125
$ws.Range("I140").Value = 6880.125
$ws.Range("J140").Value = 0
$ws.Range("K140").Value = 20640.375
$ws.Range("L140").Value = 0
$ws.Range("M140").Value = -15460.375
$ws.Range("N140").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# Row 102
$ws.Range("H102").Value = 6230.2856
$ws.Range("I102").Value = 2012
$ws.Range("K102").Value = 2012
$ws.Range("M102").Value = -390

# Row 122
$ws.Range("H122").Value = 2402459.8
$ws.Range("I122").Value = 3413047.8
$ws.Range("J122").Value = 2313.5
$ws.Range("K122").Value = 10239143.4
$ws.Range("L122").Value = 6940.5
$ws.Range("M122").Value = -10236693.4
$ws.Range("N122").Value = -11840.5

# Row 132
$ws.Range("H132").Value = 4166.4614
$ws.Range("I132").Value = 5348.75
$ws.Range("J132").Value = 3641
$ws.Range("K132").Value = 16046.25
$ws.Range("L132").Value = 10923
$ws.Range("M132").Value = -13516.25
$ws.Range("N132").Value = -15983

$ws = $wb.Worksheets.Item("LTW")
# Row 132
$ws.Range("H132").Value = 15482110
$ws.Range("I132").Value = 24081492
$ws.Range("J132").Value = 3223.8
$ws.Range("K132").Value = 72244476
$ws.Range("L132").Value = 9671.400000000001
$ws.Range("M132").Value = -72241946
$ws.Range("N132").Value = -14731.4

$ws = $wb.Worksheets.Item("WVR")
# Row 126
$ws.Range("H126").Value = 1477.9412
$ws.Range("I126").Value = 1102.5
$ws.Range("J126").Value = 2014.2858
$ws.Range("K126").Value = 3307.5
$ws.Range("L126").Value = 6042.857400000001
$ws.Range("M126").Value = -837.5
$ws.Range("N126").Value = -10982.8574

# Row 132
$ws.Range("H132").Value = 2073.5356
$ws.Range("I132").Value = 1343.5555
$ws.Range("J132").Value = 2419.3157
$ws.Range("K132").Value = 4030.6665
$ws.Range("L132").Value = 7257.9471
$ws.Range("M132").Value = -1500.6665
$ws.Range("N132").Value = -12317.9471
